$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "4th Line" - reorder the duty names, keep the same 4-person list
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("4th Line")
$ws1.Activate()
$ws1.Range("A1").Value = "Pavlo Shtefan"
$ws1.Range("A2").Value = "Olena Mikheyeva"
$ws1.Range("A3").Value = "Dmytro Latyshko"
$ws1.Range("A4").Value = "Eugene Zinchenko"
$null = $ws1.Range("A4").Select()

# ---------------------------------------------------------------------------
# Sheet "Demo" - reorder the 6-person list
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Demo")
$ws2.Range("A1").Value = "Eugene Zinchenko"
$ws2.Range("A2").Value = "Dmytro Latyshko"
$ws2.Range("A3").Value = "Olena Mikheyeva"
$ws2.Range("A4").Value = "Pavlo Shtefan"
$ws2.Range("A5").Value = "Andrii Vanikhin"
$ws2.Range("A6").Value = "Maria Donnik"

# ---------------------------------------------------------------------------
# Sheet "Night Tests" - fill the table with 5 more duty rotations (5 blocks
# of 5 people, separated by a blank row), widen column A to fit the names
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Night Tests")
$ws3.Activate()

$ws3.Cells.Item(1,1).Value = "Pavlo Shtefan"
$ws3.Cells.Item(2,1).Value = "Olena Mikheyeva"
$ws3.Cells.Item(3,1).Value = "Dmytro Latyshko"
$ws3.Cells.Item(4,1).Value = "Eugene Zinchenko"
$ws3.Cells.Item(5,1).Value = "Andrii Vanikhin"

$ws3.Cells.Item(7,1).Value = "Eugene Zinchenko"
$ws3.Cells.Item(8,1).Value = "Andrii Vanikhin"
$ws3.Cells.Item(9,1).Value = "Olena Mikheyeva"
$ws3.Cells.Item(10,1).Value = "Pavlo Shtefan"
$ws3.Cells.Item(11,1).Value = "Dmytro Latyshko"

$ws3.Cells.Item(13,1).Value = "Dmytro Latyshko"
$ws3.Cells.Item(14,1).Value = "Pavlo Shtefan"
$ws3.Cells.Item(15,1).Value = "Andrii Vanikhin"
$ws3.Cells.Item(16,1).Value = "Olena Mikheyeva"
$ws3.Cells.Item(17,1).Value = "Eugene Zinchenko"

$ws3.Cells.Item(19,1).Value = "Olena Mikheyeva"
$ws3.Cells.Item(20,1).Value = "Dmytro Latyshko"
$ws3.Cells.Item(21,1).Value = "Eugene Zinchenko"
$ws3.Cells.Item(22,1).Value = "Andrii Vanikhin"
$ws3.Cells.Item(23,1).Value = "Pavlo Shtefan"

$ws3.Cells.Item(25,1).Value = "Andrii Vanikhin"
$ws3.Cells.Item(26,1).Value = "Eugene Zinchenko"
$ws3.Cells.Item(27,1).Value = "Pavlo Shtefan"
$ws3.Cells.Item(28,1).Value = "Dmytro Latyshko"
$ws3.Cells.Item(29,1).Value = "Olena Mikheyeva"

$ws3.Columns.Item(1).ColumnWidth = 16.14
$null = $ws3.Range("A1:A29").Select()
